$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# This document has "different first page" headers/footers enabled, so each
# section exposes two relevant header/footer stories:
#   Item(1) = wdHeaderFooterPrimary   (the "default" header/footer)
#   Item(2) = wdHeaderFooterFirstPage (the "first page" header/footer)
# Both footers contain the Pearson logo picture and both headers contain the
# BTEC logo picture. Renaming an inline picture requires round-tripping it
# through a floating Shape (InlineShape has no writable Name of its own).
function Rename-InlineLogo($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    [void]$shp.ConvertToInlineShape()
}

# Pearson logo pictures (footers): image1.png -> image2.png
Rename-InlineLogo $sec.Footers.Item(1).Range "image2.png"
Rename-InlineLogo $sec.Footers.Item(2).Range "image2.png"

# BTEC logo pictures (headers): image2.jpg -> image1.jpg
Rename-InlineLogo $sec.Headers.Item(1).Range "image1.jpg"
Rename-InlineLogo $sec.Headers.Item(2).Range "image1.jpg"
